$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ch")

$ws.Range("C10").Value = "Delete image tag"
$ws.Range("C14").Value = "Close"
$ws.Range("C2").Value = "Tags"
$ws.Range("C13").Value = "Confirm"
$ws.Range("C15").Value = "Images deleted'"
